# forensic_cases.xlsx update
# - add a new "lab_results" worksheet
# - add created_at columns to users / labs / samples
# - rename several id-ish columns to friendlier names across cases / samples / custody_events
# - seed an initial admin user row in users

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

function Copy-HeaderStyle($srcRange, $dstRange) {
    $srcRange.Copy() | Out-Null
    $dstRange.PasteSpecial($xlPasteFormats) | Out-Null
}

# ---------------------------------------------------------------------------
# users: add "created_at" header, and seed the initial admin user (row 2)
# ---------------------------------------------------------------------------
$wsUsers = $wb.Worksheets.Item("users")

Copy-HeaderStyle $wsUsers.Range("F1") $wsUsers.Range("G1")
$wsUsers.Range("G1").Value = "created_at"

$wsUsers.Range("A2").Value = 1
$wsUsers.Range("B2").Value = "admin@fasttrack.local"
$wsUsers.Range("C2").Value = "Admin"
$wsUsers.Range("D2").Value = "admin"
$wsUsers.Range("E2").Value = "scrypt:32768:8:1`$2jDA8TMQ61QO1h39`$05fa170edb747deedcfd99b0a23d8c004425219406ccb6a0309c5b9dc6df63ac0b5ac98dd4d238fed4523aa85a33029efd8496c727785536e3b280b0d3880a0c"
# api_token stays blank for the seeded admin row, but touch it so the row's
# used-range still stretches out to column F like the source data.
$wsUsers.Range("F2").Value = ""
$wsUsers.Range("G2").Value = "2025-08-31T10:56:41.919440"

# ---------------------------------------------------------------------------
# labs: add "created_at" header
# ---------------------------------------------------------------------------
$wsLabs = $wb.Worksheets.Item("labs")

Copy-HeaderStyle $wsLabs.Range("D1") $wsLabs.Range("E1")
$wsLabs.Range("E1").Value = "created_at"

# ---------------------------------------------------------------------------
# cases: rename created_by_id -> created_by, lab_id -> lab_assigned
# ---------------------------------------------------------------------------
$wsCases = $wb.Worksheets.Item("cases")

$wsCases.Range("H1").Value = "created_by"
$wsCases.Range("I1").Value = "lab_assigned"

# ---------------------------------------------------------------------------
# samples: rename case_id -> case_number, add "created_at" header
# ---------------------------------------------------------------------------
$wsSamples = $wb.Worksheets.Item("samples")

$wsSamples.Range("B1").Value = "case_number"

Copy-HeaderStyle $wsSamples.Range("E1") $wsSamples.Range("F1")
$wsSamples.Range("F1").Value = "created_at"

# ---------------------------------------------------------------------------
# custody_events: rename case_id -> case_number, sample_id -> sample_code,
# actor_id -> actor
# ---------------------------------------------------------------------------
$wsCustody = $wb.Worksheets.Item("custody_events")

$wsCustody.Range("B1").Value = "case_number"
$wsCustody.Range("C1").Value = "sample_code"
$wsCustody.Range("D1").Value = "actor"

# ---------------------------------------------------------------------------
# lab_results: brand new worksheet, appended after custody_events
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsLabResults = $wb.Worksheets.Add($null, $lastSheet)
$wsLabResults.Name = "lab_results"

$headers = @("id", "case_number", "sample_code", "lab_user", "result_summary", "result_file", "created_at")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $cell = $wsLabResults.Cells.Item(1, $col)
    Copy-HeaderStyle $wsUsers.Range("A1") $cell
    $cell.Value = $headers[$i]
}

$excel.CutCopyMode = $false

Write-Host "lab_results sheet created and headers updated"
